# "Generate Report for Archive"
# - The handoff status text moves from "Ready for handoff" to "In Translation"
#   everywhere it appears (Overview!E2:E3/F2:F3 and the zh-cn / de-de sheets'
#   Status column C2:C3).
# - Because the displayed text got shorter, the Status columns that were
#   sized to fit it shrink too (Overview columns E:F, and column C on the
#   zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

# --- Update the status text wherever it is used ---------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Shrink the columns that held the (now shorter) status text -----------
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
